$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $Val)
    $r = $Sheet.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "69.541.00"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "3.788.88"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue $ws "D5" "613.98"
$ws.Range("E5").Value = "  -0.92%  "
Set-TextValue $ws "D6" "176.87"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "3.785.77"
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("E12").Value = "  -0.99%  "
Set-TextValue $ws "D13" "39.81"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "4.418.19"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "3.786.13"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "69.606.09"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  -3.56%  "
Set-TextValue $ws "D20" "16.62"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  +0.55%  "
Set-TextValue $ws "D22" "9.66"
$ws.Range("E22").Value = "  +1.40%  "
Set-TextValue $ws "D23" "0.735"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E24").Value = "  -1.29%  "
Set-TextValue $ws "D25" "86.26"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  +4.31%  "
Set-TextValue $ws "D27" "12.81"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("E28").Value = "  -5.17%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("E31").Value = "  +0.00%  "
Set-TextValue $ws "D32" "8.10"
$ws.Range("E32").Value = "  +2.52%  "
Set-TextValue $ws "D33" "31.33"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("E37").Value = "  -1.28%  "
Set-TextValue $ws "D38" "0.142"
$ws.Range("E38").Value = "  +7.58%  "
Set-TextValue $ws "D39" "484.58"
$ws.Range("E39").Value = "  +13.23%  "
Set-TextValue $ws "D40" "0.339"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D42" "2.99"
$ws.Range("E42").Value = "  +4.61%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D43" "49.77"
$ws.Range("E43").Value = "  -0.79%  "
Set-TextValue $ws "D44" "44.21"
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").Value = "2.941.16"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -0.27%  "
Set-TextValue $ws "D48" "27.41"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D49" "1.00"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D50" "138.91"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("E51").Value = "  -1.51%  "
